$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2009年 data row (row 2). This shifts the 2010年 row (row 3)
# up to become the new row 2, matching the diff.
$ws.Range("A2:U2").EntireRow.Delete()
